# Updated cryptos list on Mon Nov 18 22:44:02 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and swaps row 51 from VeChain to OKB.
#
# Column D values are prefixed with a leading apostrophe so Excel stores
# them as literal text (matching the workbook's existing inline-string
# cells) instead of auto-coercing them into numbers, which would silently
# drop meaningful trailing/leading zeros (e.g. "240.40" -> 240.4) or mangle
# the multi-dot "thousands" style prices (e.g. "91.411.90").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'91.411.90"
$ws.Range('E2').Value = '  +2.23%  '
$ws.Range('D3').Value = "'3.162.23"
$ws.Range('E3').Value = '  +3.14%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'240.40"
$ws.Range('E5').Value = '  +2.73%  '
$ws.Range('D6').Value = "'620.69"
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('E7').Value = '  +6.59%  '
$ws.Range('D8').Value = "'0.376"
$ws.Range('E8').Value = '  +4.24%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('D10').Value = "'3.159.29"
$ws.Range('E10').Value = '  +2.91%  '
$ws.Range('D11').Value = "'0.747"
$ws.Range('E11').Value = '  +5.92%  '
$ws.Range('E12').Value = '  +2.41%  '
$ws.Range('D13').Value = "'0.0000248"
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('E14').Value = '  +1.93%  '
$ws.Range('E15').Value = '  +5.05%  '
$ws.Range('D16').Value = "'91.211.24"
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').Value = "'3.746.16"
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').Value = "'3.174.84"
$ws.Range('E18').Value = '  +3.19%  '
$ws.Range('D19').Value = "'3.77"
$ws.Range('E19').Value = '  +1.78%  '
$ws.Range('D20').Value = "'15.26"
$ws.Range('E20').Value = '  +11.55%  '
$ws.Range('D21').Value = "'6.06"
$ws.Range('E21').Value = '  +12.88%  '
$ws.Range('D22').Value = "'457.00"
$ws.Range('E22').Value = '  +6.50%  '
$ws.Range('E23').Value = '  -4.10%  '
$ws.Range('E24').Value = '  +6.86%  '
$ws.Range('E25').Value = '  +9.19%  '
$ws.Range('D26').Value = "'89.21"
$ws.Range('E26').Value = '  +2.60%  '
$ws.Range('D27').Value = "'12.10"
$ws.Range('E27').Value = '  +4.39%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = "'0.132"
$ws.Range('E30').Value = '  +47.38%  '
$ws.Range('D31').Value = "'0.235"
$ws.Range('E31').Value = '  +18.74%  '
$ws.Range('D32').Value = "'0.172"
$ws.Range('E32').Value = '  +10.36%  '
$ws.Range('D33').Value = "'9.44"
$ws.Range('E33').Value = '  +5.23%  '
$ws.Range('E34').Value = '  +15.66%  '
$ws.Range('E35').Value = '  -10.64%  '
$ws.Range('D36').Value = "'7.71"
$ws.Range('E36').Value = '  +10.65%  '
$ws.Range('D37').Value = "'26.57"
$ws.Range('E37').Value = '  +4.25%  '
$ws.Range('D38').Value = "'515.43"
$ws.Range('E38').Value = '  +5.64%  '
$ws.Range('D39').Value = "'1.97"
$ws.Range('E39').Value = '  +5.30%  '
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('D41').Value = "'3.94"
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').Value = "'0.453"
$ws.Range('E42').Value = '  +14.81%  '
$ws.Range('D43').Value = "'3.48"
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('D44').Value = "'22.15"
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('E46').Value = '  +7.41%  '
$ws.Range('E47').Value = '  +6.52%  '
$ws.Range('D48').Value = "'158.36"
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('E49').Value = '  +7.13%  '
$ws.Range('D50').Value = "'4.50"
$ws.Range('E50').Value = '  +4.72%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = "'44.01"
$ws.Range('E51').Value = '  -0.40%  '
